$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.830.41'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +7.91%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.808.95'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +5.00%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.19%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '248.06'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.32%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("D6").Style = "Normal"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4958'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.31%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2776'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +7.85%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06413'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.75%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.809.20'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +4.95%  '

$ws.Range("E11").Value = '  +5.52%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07067'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.74%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.6466'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +7.22%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '84.03'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +9.42%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.686'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +5.14%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '28.813.52'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +8.53%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.000'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.16%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000007338'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +3.01%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.0000'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.14%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.24'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +7.89%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.047.83'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +5.02%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.569'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.73%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.899'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.94%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.340'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +5.76%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '142.35'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.61%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '129.47'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +21.77%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.40'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +7.77%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.881'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +6.39%  '

$ws.Range("E29").Value = '  +3.55%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.131'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.21%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08340'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +5.49%  '

$ws.Range("E32").Value = '  +3.75%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04954'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +10.89%  '

$ws.Range("E34").Value = '  +9.11%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.720'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +4.76%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6701'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +8.83%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.261'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +12.82%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.726'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +11.51%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.9569'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.34%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.081'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +8.62%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.01588'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +6.75%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9998'
$ws.Range("D42").Style = "Normal"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.4072'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +6.58%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '7.153'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +5.69%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05514'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.64%  '

$ws.Range("E48").Value = '  +3.61%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '31.56'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +5.01%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.3619'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +8.36%  '

$ws.Range("E51").Value = '  +6.09%  '
